# Change rule name in B11 from "R40" to "1".
# The new value must remain stored as text (shared string), matching the
# original cell's data type, and the cell's existing style/formatting
# (s="23") must be preserved unchanged.
#
# A plain `$ws.Range("B11").Value = "1"` would be auto-coerced to a number
# (since "1" looks numeric and the cell's number format is General), and
# forcing text via NumberFormat="@" on B11 directly would permanently change
# its style index. So we stage the text value on a scratch cell, copy it
# over (value+type only), then restore B11's original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("ZZ9999")

# Back up B11's current formatting onto the scratch cell.
$ws.Range("B11").Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats

# Force the new value to be stored as text, not a number.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"

# Restore B11's original formatting (style), now that the text value is set.
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats

# Clean up the scratch cell so nothing else changes in the sheet.
$scratch.Clear()

$excel.CutCopyMode = $false
